$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1619
$ws.Range("E2").Value = 207
$ws.Range("F2").Value = 223
$ws.Range("G2").Value = 249
$ws.Range("H2").Value = 169
$ws.Range("I2").Value = 78
$ws.Range("J2").Value = 91
$ws.Range("K2").Value = 2396
$ws.Range("L2").Value = 349
$ws.Range("M2").Value = 2047
$ws.Range("N2").Value = 1228
$ws.Range("O2").Value = 818
$ws.Range("P2").Value = 75
$ws.Range("Q2").Value = 117
$ws.Range("R2").Value = -70
$ws.Range("S2").Value = -57
$ws.Range("T2").Value = 47
$ws.Range("U2").Value = 70
$ws.Range("V2").Value = 33
$ws.Range("W2").Value = 12.8
$ws.Range("X2").Value = 10.45
$ws.Range("Y2").Value = 6.46
$ws.Range("Z2").Value = 7.21
$ws.Range("AA2").Value = 17.06
$ws.Range("AB2").Value = 1694.8
$ws.Range("AC2").Value = 523
$ws.Range("AD2").Value = 8.91
$ws.Range("AE2").Value = 9063
$ws.Range("AF2").Value = 0.51
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 4.29
$ws.Range("AI2").Value = 34.72
$ws.Range("AJ2").Value = 14934008

# Row 3
$ws.Range("D3").Value = 1851
$ws.Range("E3").Value = 235
$ws.Range("F3").Value = 260
$ws.Range("G3").Value = 333
$ws.Range("H3").Value = 248
$ws.Range("I3").Value = 112
$ws.Range("J3").Value = 136
$ws.Range("K3").Value = 2603
$ws.Range("L3").Value = 441
$ws.Range("M3").Value = 2162
$ws.Range("N3").Value = 1281
$ws.Range("O3").Value = 881
$ws.Range("P3").Value = 75
$ws.Range("Q3").Value = 272
$ws.Range("R3").Value = -186
$ws.Range("S3").Value = 33
$ws.Range("T3").Value = 88
$ws.Range("U3").Value = 183
$ws.Range("V3").Value = 120
$ws.Range("W3").Value = 12.72
$ws.Range("X3").Value = 13.39
$ws.Range("Y3").Value = 8.9
$ws.Range("Z3").Value = 9.92
$ws.Range("AA3").Value = 20.38
$ws.Range("AB3").Value = 1805.69
$ws.Range("AC3").Value = 747
$ws.Range("AD3").Value = 7.77
$ws.Range("AE3").Value = 9455
$ws.Range("AF3").Value = 0.61
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 3.44
$ws.Range("AI3").Value = 24.28
$ws.Range("AJ3").Value = 14934008

# Row 4
$ws.Range("D4").Value = 1814
$ws.Range("E4").Value = 136
$ws.Range("F4").Value = 136
$ws.Range("G4").Value = 205
$ws.Range("H4").Value = 152
$ws.Range("I4").Value = 57
$ws.Range("J4").Value = 95
$ws.Range("K4").Value = 2788
$ws.Range("L4").Value = 482
$ws.Range("M4").Value = 2307
$ws.Range("N4").Value = 1327
$ws.Range("O4").Value = 980
$ws.Range("P4").Value = 75
$ws.Range("Q4").Value = 175
$ws.Range("R4").Value = -194
$ws.Range("S4").Value = -48
$ws.Range("T4").Value = 68
$ws.Range("U4").Value = 108
$ws.Range("V4").Value = 123
$ws.Range("W4").Value = 7.51
$ws.Range("X4").Value = 8.359999999999999
$ws.Range("Y4").Value = 4.35
$ws.Range("Z4").Value = 5.62
$ws.Range("AA4").Value = 20.89
$ws.Range("AB4").Value = 1857.66
$ws.Range("AC4").Value = 380
$ws.Range("AD4").Value = 13.75
$ws.Range("AE4").Value = 9790
$ws.Range("AF4").Value = 0.53
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 3.83
$ws.Range("AI4").Value = 47.81
$ws.Range("AJ4").Value = 14934008

# Row 5
$ws.Range("D5").Value = 1698
$ws.Range("E5").Value = 104
$ws.Range("F5").Value = 104
$ws.Range("G5").Value = 282
$ws.Range("H5").Value = 193
$ws.Range("I5").Value = 153
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 2917
$ws.Range("L5").Value = 476
$ws.Range("M5").Value = 2441
$ws.Range("N5").Value = 1451
$ws.Range("O5").Value = 990
$ws.Range("P5").Value = 75
$ws.Range("Q5").Value = 78
$ws.Range("R5").Value = 23
$ws.Range("S5").Value = -49
$ws.Range("T5").Value = 48
$ws.Range("U5").Value = 30
$ws.Range("V5").Value = 125
$ws.Range("W5").Value = 6.15
$ws.Range("X5").Value = 11.38
$ws.Range("Y5").Value = 11.02
$ws.Range("Z5").Value = 6.77
$ws.Range("AA5").Value = 19.51
$ws.Range("AB5").Value = 2037.52
$ws.Range("AC5").Value = 1025
$ws.Range("AD5").Value = 4.68
$ws.Range("AE5").Value = 10707
$ws.Range("AF5").Value = 0.45
$ws.Range("AG5").Value = 200
$ws.Range("AH5").Value = 4.17
$ws.Range("AI5").Value = 17.7
$ws.Range("AJ5").Value = 14934008

# Row 6
$ws.Range("D6").Value = 1607
$ws.Range("E6").Value = 18
$ws.Range("F6").Value = 18
$ws.Range("G6").Value = 107
$ws.Range("H6").Value = 62
$ws.Range("I6").Value = 31
$ws.Range("K6").Value = 2895
$ws.Range("L6").Value = 456
$ws.Range("M6").Value = 2439
$ws.Range("N6").Value = 1463
$ws.Range("P6").Value = 75
$ws.Range("Q6").Value = 117
$ws.Range("R6").Value = -110
$ws.Range("S6").Value = -35
$ws.Range("T6").Value = 73
$ws.Range("U6").Value = 44
$ws.Range("V6").Value = 133
$ws.Range("W6").Value = 1.11
$ws.Range("X6").Value = 3.83
$ws.Range("Y6").Value = 2.11
$ws.Range("Z6").Value = 2.12
$ws.Range("AA6").Value = 18.7
$ws.Range("AB6").Value = 2019.91
$ws.Range("AC6").Value = 206
$ws.Range("AD6").Value = 16.64
$ws.Range("AE6").Value = 10796
$ws.Range("AF6").Value = 0.32
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 4.38
$ws.Range("AI6").Value = 66.14
$ws.Range("AJ6").Value = 14934008

# Clear rows 7-9 (D:AI) - data no longer reported
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()